$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 187.41176
$ws.Range("I2").Value = 186.35715
$ws.Range("J2").Value = 192.33333
$ws.Range("K2").Value = 186.35715
$ws.Range("L2").Value = 192.33333
$ws.Range("M2").Value = -73.35714999999999
$ws.Range("N2").Value = -418.33333
$ws.Range("H6").Value = 270.875
$ws.Range("I6").Value = 61.166668
$ws.Range("K6").Value = 183.500004
$ws.Range("M6").Value = -71.50000399999999
$ws.Range("H32").Value = 14292642
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H40").Value = 3561.5667
$ws.Range("I40").Value = 2497.6875
$ws.Range("J40").Value = 4777.4287
$ws.Range("K40").Value = 2497.6875
$ws.Range("L40").Value = 4777.4287
$ws.Range("M40").Value = -2322.6875
$ws.Range("N40").Value = -5127.4287
$ws.Range("H53").Value = 611.2
$ws.Range("I53").Value = 216.66667
$ws.Range("J53").Value = 780.2857
$ws.Range("K53").Value = 216.66667
$ws.Range("L53").Value = 780.2857
$ws.Range("M53").Value = 420.33333
$ws.Range("N53").Value = -2054.2857
$ws.Range("H86").Value = 3762701.5
$ws.Range("I86").Value = 3511.2856
$ws.Range("J86").Value = 7521891.5
$ws.Range("K86").Value = 3511.2856
$ws.Range("L86").Value = 7521891.5
$ws.Range("M86").Value = -2388.2856
$ws.Range("N86").Value = -7524137.5
$ws.Range("H89").Value = 3762701.5
$ws.Range("I89").Value = 3511.2856
$ws.Range("J89").Value = 7521891.5
$ws.Range("K89").Value = 17556.428
$ws.Range("L89").Value = 37609457.5
$ws.Range("M89").Value = -11940.428
$ws.Range("N89").Value = -37620689.5
$ws.Range("H111").Value = 13336570
$ws.Range("I111").Value = 18184114
$ws.Range("J111").Value = 5821.75
$ws.Range("K111").Value = 54552342
$ws.Range("L111").Value = 17465.25
$ws.Range("M111").Value = -54549275
$ws.Range("N111").Value = -23599.25
$ws.Range("H113").Value = 10327
$ws.Range("J113").Value = 10492.857
$ws.Range("L113").Value = 10492.857
$ws.Range("N113").Value = -17000.857
$ws.Range("H138").Value = 1303.0781
$ws.Range("J138").Value = 4515.6665
$ws.Range("L138").Value = 13546.9995
$ws.Range("N138").Value = -23826.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22733558
$ws.Range("I32").Value = 23815918
$ws.Range("J32").Value = 4000
$ws.Range("K32").Value = 23815918
$ws.Range("L32").Value = 4000
$ws.Range("M32").Value = -23815631
$ws.Range("N32").Value = -4574
$ws.Range("H61").Value = 3125.4614
$ws.Range("I61").Value = 1853.5555
$ws.Range("J61").Value = 5987.25
$ws.Range("K61").Value = 1853.5555
$ws.Range("L61").Value = 5987.25
$ws.Range("M61").Value = -1641.5555
$ws.Range("N61").Value = -6411.25
$ws.Range("H74").Value = 1377.3334
$ws.Range("I74").Value = 1016.65
$ws.Range("K74").Value = 1016.65
$ws.Range("M74").Value = -142.65
$ws.Range("H77").Value = 1377.3334
$ws.Range("I77").Value = 1016.65
$ws.Range("K77").Value = 5083.25
$ws.Range("M77").Value = -715.25
$ws.Range("H102").Value = 9805104
$ws.Range("I102").Value = 1244.7742
$ws.Range("K102").Value = 1244.7742
$ws.Range("M102").Value = 377.2257999999999
$ws.Range("H122").Value = 1589.5588
$ws.Range("I122").Value = 1363.6207
$ws.Range("K122").Value = 4090.8621
$ws.Range("M122").Value = -1640.8621
$ws.Range("H132").Value = 1318.5834
$ws.Range("I132").Value = 1311.8667
$ws.Range("K132").Value = 3935.6001
$ws.Range("M132").Value = -1405.6001
$ws.Range("H136").Value = 3125.4614
$ws.Range("I136").Value = 1853.5555
$ws.Range("J136").Value = 5987.25
$ws.Range("K136").Value = 5560.666499999999
$ws.Range("L136").Value = 17961.75
$ws.Range("M136").Value = -3010.666499999999
$ws.Range("N136").Value = -23061.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2476.054
$ws.Range("J86").Value = 3799.158
$ws.Range("L86").Value = 3799.158
$ws.Range("N86").Value = -6045.157999999999
$ws.Range("H89").Value = 2476.054
$ws.Range("J89").Value = 3799.158
$ws.Range("L89").Value = 18995.79
$ws.Range("N89").Value = -30227.79
$ws.Range("H99").Value = 1636.4524
$ws.Range("I99").Value = 1486.5454
$ws.Range("K99").Value = 1486.5454
$ws.Range("M99").Value = 11.45460000000003
$ws.Range("H134").Value = 2215.0688
$ws.Range("I134").Value = 1693.7556
$ws.Range("K134").Value = 5081.266799999999
$ws.Range("M134").Value = -2546.266799999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3234.4695
$ws.Range("I134").Value = 2283.139
$ws.Range("K134").Value = 6849.417
$ws.Range("M134").Value = -4314.417

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 136.44444
$ws.Range("I12").Value = 156.66667
$ws.Range("J12").Value = 132.4
$ws.Range("K12").Value = 470.00001
$ws.Range("L12").Value = 397.2
$ws.Range("M12").Value = -297.00001
$ws.Range("N12").Value = -743.2
$ws.Range("H109").Value = 1304.6471
$ws.Range("H131").Value = 3030.5
$ws.Range("I131").Value = 1704.3846
$ws.Range("J131").Value = 4356.615
$ws.Range("K131").Value = 5113.1538
$ws.Range("L131").Value = 13069.845
$ws.Range("M131").Value = -73.15380000000005
$ws.Range("N131").Value = -23149.845
$ws.Range("H132").Value = 1060.4
$ws.Range("I132").Value = 1060.4
$ws.Range("K132").Value = 9543.6
$ws.Range("M132").Value = -7013.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 2500
$ws.Range("I10").Value = 2500
$ws.Range("K10").Value = 2500
$ws.Range("M10").Value = -2331
$ws.Range("H11").Value = 2603555.5
$ws.Range("I11").Value = 2051333.4
$ws.Range("K11").Value = 2051333.4
$ws.Range("M11").Value = -2051194.4
$ws.Range("H70").Value = 78138.25
$ws.Range("I70").Value = 164216.58
$ws.Range("K70").Value = 164216.58
$ws.Range("M70").Value = -163946.58
$ws.Range("H73").Value = 78138.25
$ws.Range("I73").Value = 164216.58
$ws.Range("K73").Value = 164216.58
$ws.Range("M73").Value = -163280.58
$ws.Range("H97").Value = 8333722.5
$ws.Range("I97").Value = 418.58334
$ws.Range("K97").Value = 418.58334
$ws.Range("M97").Value = 77.41665999999998
$ws.Range("H102").Value = 3333.9092
$ws.Range("I102").Value = 1767.5
$ws.Range("K102").Value = 1767.5
$ws.Range("M102").Value = -145.5
$ws.Range("H107").Value = 1151
$ws.Range("I107").Value = 1548.2858
$ws.Range("K107").Value = 1548.2858
$ws.Range("M107").Value = 371.7141999999999
$ws.Range("H113").Value = 3178.7878
$ws.Range("J113").Value = 8540.25
$ws.Range("L113").Value = 8540.25
$ws.Range("N113").Value = -12880.25
$ws.Range("H132").Value = 2918
$ws.Range("I132").Value = 2754.6
$ws.Range("K132").Value = 8263.799999999999
$ws.Range("M132").Value = -5733.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 451.5
$ws.Range("I16").Value = 373.14285
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 373.14285
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -203.14285
$ws.Range("N16").Value = -1340
$ws.Range("H36").Value = 65874.25
$ws.Range("J36").Value = 65874.25
$ws.Range("L36").Value = 65874.25
$ws.Range("N36").Value = -66998.25
$ws.Range("H40").Value = 9835.357
$ws.Range("I40").Value = 11218.637
$ws.Range("K40").Value = 11218.637
$ws.Range("M40").Value = -11082.637
$ws.Range("H61").Value = 766
$ws.Range("I61").Value = 673.75
$ws.Range("K61").Value = 673.75
$ws.Range("M61").Value = -471.75
$ws.Range("H68").Value = 1685.5555
$ws.Range("I68").Value = 1685.5555
$ws.Range("K68").Value = 1685.5555
$ws.Range("M68").Value = -936.5554999999999
$ws.Range("H71").Value = 1685.5555
$ws.Range("I71").Value = 1685.5555
$ws.Range("K71").Value = 8427.7775
$ws.Range("M71").Value = -4683.7775
$ws.Range("H93").Value = 11907912
$ws.Range("I93").Value = 3494
$ws.Range("K93").Value = 3494
$ws.Range("M93").Value = -2246
$ws.Range("H113").Value = 766
$ws.Range("I113").Value = 673.75
$ws.Range("K113").Value = 673.75
$ws.Range("M113").Value = 1496.25
$ws.Range("H122").Value = 6036.278
$ws.Range("I122").Value = 2900.5
$ws.Range("K122").Value = 8701.5
$ws.Range("M122").Value = -6251.5
$ws.Range("H132").Value = 3482.0364
$ws.Range("I132").Value = 2411.7446
$ws.Range("K132").Value = 7235.2338
$ws.Range("M132").Value = -4705.2338
$ws.Range("H136").Value = 5071.4614
$ws.Range("I136").Value = 3742.8235
$ws.Range("J136").Value = 7581.1113
$ws.Range("K136").Value = 11228.4705
$ws.Range("L136").Value = 22743.3339
$ws.Range("M136").Value = -8678.470499999999
$ws.Range("N136").Value = -27843.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 594759.9
$ws.Range("I122").Value = 3180.2144
$ws.Range("K122").Value = 9540.643199999999
$ws.Range("M122").Value = -7090.643199999999
$ws.Range("H132").Value = 2846.4783
$ws.Range("I132").Value = 2551.0527
$ws.Range("J132").Value = 4249.75
$ws.Range("K132").Value = 7653.158100000001
$ws.Range("L132").Value = 12749.25
$ws.Range("M132").Value = -5123.158100000001
$ws.Range("N132").Value = -17809.25
$ws.Range("H136").Value = 964.0784
$ws.Range("I136").Value = 522.87177
$ws.Range("J136").Value = 2398
$ws.Range("K136").Value = 1568.61531
$ws.Range("L136").Value = 7194
$ws.Range("M136").Value = 981.3846900000001
$ws.Range("N136").Value = -12294
